$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New ActualRate (column E) values for rows 2-31. These are dollar-amount
# strings stored as plain text (not currency numbers) in the source data, so
# we round-trip them through a scratch cell's text formula + copy/paste
# (values only) to avoid Excel's "looks like currency" auto-conversion that a
# direct .Value assignment of a string like "$18.40" would trigger.
$newRates = @(
    "$18.40",
    "$19.73",
    "$27.61",
    "$38.85",
    "$41.36",
    "$63.39",
    "$74.80",
    "$210.32",
    "$6.00",
    "$6.00",
    "$24.00",
    "$36.00",
    "$6.00",
    "$63.00",
    "$6.18",
    "$13.55",
    "$18.91",
    "$29.42",
    "$35.13",
    "$48.33",
    "$37.82",
    "$48.33",
    "$56.73",
    "$79.85",
    "$105.06",
    "$105.06",
    "$148.16",
    "$287.20",
    "$476.40",
    "$273.68"
)

$firstRow = 2
$scratch = $ws.Cells.Item(100, 1)

for ($i = 0; $i -lt $newRates.Length; $i++) {
    $row = $firstRow + $i
    $scratch.Formula = '="' + $newRates[$i] + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, 5).PasteSpecial(-4163)
}
$scratch.Clear()

# Column F (Result) simply flips to FAIL for every one of these rows.
$ws.Range("F2:F31").Value = "FAIL"
